{"js": "// Auto-generated edit script.\n// Replaces the date paragraph text and the 100 arithmetic-expression table\n// cells in document order (paragraph 0 = date paragraph above the table;\n// paragraphs 1..100 = the 20x5 table cells, row-major, matching the\n// body.paragraphs enumeration order).\nconst OLD_VALUES = [\"2025-04-24 Thursday\", \"22-2=\", \"78-6=\", \"5+51=\", \"74-46=\", \"56-48=\", \"39+16=\", \"36-15=\", \"88-52=\", \"26+11=\", \"37+8=\", \"33+11=\", \"23+76=\", \"22+73=\", \"33+63=\", \"93-77=\", \"3-3=\", \"70-1=\", \"24+55=\", \"6+75=\", \"40+25=\", \"11-6=\", \"0+40=\", \"36+55=\", \"14+2=\", \"9+36=\", \"89-14=\", \"85+3=\", \"6+0=\", \"76-59=\", \"3+65=\", \"59-17=\", \"84-73=\", \"83-39=\", \"80+8=\", \"12+32=\", \"45-24=\", \"73-48=\", \"47-27=\", \"82+13=\", \"82-0=\", \"59-2=\", \"11+2=\", \"6+13=\", \"57-49=\", \"75-31=\", \"11+21=\", \"79+12=\", \"83-7=\", \"97-9=\", \"82-5=\", \"80-34=\", \"51+13=\", \"30+67=\", \"4+19=\", \"86-51=\", \"29+54=\", \"2+20=\", \"14+72=\", \"55-42=\", \"41+51=\", \"45-40=\", \"41-14=\", \"94-51=\", \"44-16=\", \"45-44=\", \"2+31=\", \"74-58=\", \"19+22=\", \"47+23=\", \"51+3=\", \"38+35=\", \"63-15=\", \"52+20=\", \"41+18=\", \"65+19=\", \"27+54=\", \"51+30=\", \"88-64=\", \"59-19=\", \"47+41=\", \"30-4=\", \"24-2=\", \"97-60=\", \"89-87=\", \"97-60=\", \"17+40=\", \"72+26=\", \"37-10=\", \"45-39=\", \"46-10=\", \"29+46=\", \"24+14=\", \"8+82=\", \"54+34=\", \"58+37=\", \"64+19=\", \"4+94=\", \"51-12=\", \"36+43=\", \"58+11=\"];\nconst NEW_VALUES = [\"2025-04-25 Friday\", \"56-3=\", \"77-24=\", \"26+47=\", \"73-22=\", \"51-14=\", \"1+12=\", \"17-17=\", \"45+20=\", \"51-16=\", \"92-3=\", \"11+30=\", \"20+24=\", \"8+62=\", \"19+33=\", \"61+26=\", \"95-50=\", \"85-0=\", \"74-35=\", \"24+12=\", \"21+10=\", \"93-53=\", \"7+21=\", \"34+27=\", \"10+42=\", \"81-14=\", \"42+27=\", \"28+18=\", \"25+64=\", \"95-32=\", \"63-45=\", \"8+25=\", \"44+40=\", \"29-0=\", \"13+15=\", \"25+3=\", \"39+29=\", \"6+14=\", \"57-33=\", \"57-29=\", \"31-22=\", \"44+16=\", \"66-30=\", \"56-55=\", \"45+25=\", \"1+87=\", \"52-7=\", \"6+28=\", \"58-7=\", \"65-43=\", \"49+37=\", \"36-2=\", \"89-4=\", \"66-27=\", \"57-22=\", \"97-3=\", \"11+33=\", \"47+34=\", \"58-38=\", \"18+71=\", \"29+12=\", \"11+85=\", \"23+61=\", \"8+24=\", \"69-5=\", \"67-4=\", \"96-58=\", \"0+81=\", \"39-14=\", \"78-17=\", \"92-72=\", \"45+43=\", \"15-11=\", \"17+53=\", \"49+17=\", \"46-8=\", \"64-19=\", \"69-66=\", \"20+43=\", \"65+24=\", \"28+1=\", \"69-53=\", \"36+4=\", \"45+10=\", \"47-41=\", \"72-58=\", \"87-76=\", \"12+8=\", \"85+7=\", \"47-43=\", \"26-17=\", \"74+8=\", \"78-73=\", \"22+45=\", \"99-32=\", \"94-17=\", \"36+13=\", \"53-13=\", \"49+30=\", \"69-20=\", \"44+34=\"];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== OLD_VALUES.length) {\n  throw new Error(\n    \"Unexpected paragraph count: expected \" + OLD_VALUES.length +\n    \" got \" + paragraphs.items.length\n  );\n}\n\n// Load current text for every paragraph so we can sanity-check before writing.\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const expected = OLD_VALUES[i];\n  const actual = para.text;\n  if (actual !== expected) {\n    throw new Error(\n      \"Paragraph \" + i + \" text mismatch: expected \" + JSON.stringify(expected) +\n      \" got \" + JSON.stringify(actual)\n    );\n  }\n  if (NEW_VALUES[i] !== expected) {\n    para.insertText(NEW_VALUES[i], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Auto-generated edit script.\n# Updates the date line above the table and the 100 arithmetic-expression\n# cells inside the 20x5 table, addressed positionally (row-major) so the\n# two duplicate \"97-60=\" cells (row 17, columns 3 & 5) are each mapped to\n# their own distinct replacement instead of relying on ambiguous text search.\n\n$d = $word.ActiveDocument\n\n$dateOld = '2025-04-24 Thursday'\n$dateNew = '2025-04-25 Friday'\n\n$oldValues = @(\n    '22-2=',\n    '78-6=',\n    '5+51=',\n    '74-46=',\n    '56-48=',\n    '39+16=',\n    '36-15=',\n    '88-52=',\n    '26+11=',\n    '37+8=',\n    '33+11=',\n    '23+76=',\n    '22+73=',\n    '33+63=',\n    '93-77=',\n    '3-3=',\n    '70-1=',\n    '24+55=',\n    '6+75=',\n    '40+25=',\n    '11-6=',\n    '0+40=',\n    '36+55=',\n    '14+2=',\n    '9+36=',\n    '89-14=',\n    '85+3=',\n    '6+0=',\n    '76-59=',\n    '3+65=',\n    '59-17=',\n    '84-73=',\n    '83-39=',\n    '80+8=',\n    '12+32=',\n    '45-24=',\n    '73-48=',\n    '47-27=',\n    '82+13=',\n    '82-0=',\n    '59-2=',\n    '11+2=',\n    '6+13=',\n    '57-49=',\n    '75-31=',\n    '11+21=',\n    '79+12=',\n    '83-7=',\n    '97-9=',\n    '82-5=',\n    '80-34=',\n    '51+13=',\n    '30+67=',\n    '4+19=',\n    '86-51=',\n    '29+54=',\n    '2+20=',\n    '14+72=',\n    '55-42=',\n    '41+51=',\n    '45-40=',\n    '41-14=',\n    '94-51=',\n    '44-16=',\n    '45-44=',\n    '2+31=',\n    '74-58=',\n    '19+22=',\n    '47+23=',\n    '51+3=',\n    '38+35=',\n    '63-15=',\n    '52+20=',\n    '41+18=',\n    '65+19=',\n    '27+54=',\n    '51+30=',\n    '88-64=',\n    '59-19=',\n    '47+41=',\n    '30-4=',\n    '24-2=',\n    '97-60=',\n    '89-87=',\n    '97-60=',\n    '17+40=',\n    '72+26=',\n    '37-10=',\n    '45-39=',\n    '46-10=',\n    '29+46=',\n    '24+14=',\n    '8+82=',\n    '54+34=',\n    '58+37=',\n    '64+19=',\n    '4+94=',\n    '51-12=',\n    '36+43=',\n    '58+11='\n)\n$newValues = @(\n    '56-3=',\n    '77-24=',\n    '26+47=',\n    '73-22=',\n    '51-14=',\n    '1+12=',\n    '17-17=',\n    '45+20=',\n    '51-16=',\n    '92-3=',\n    '11+30=',\n    '20+24=',\n    '8+62=',\n    '19+33=',\n    '61+26=',\n    '95-50=',\n    '85-0=',\n    '74-35=',\n    '24+12=',\n    '21+10=',\n    '93-53=',\n    '7+21=',\n    '34+27=',\n    '10+42=',\n    '81-14=',\n    '42+27=',\n    '28+18=',\n    '25+64=',\n    '95-32=',\n    '63-45=',\n    '8+25=',\n    '44+40=',\n    '29-0=',\n    '13+15=',\n    '25+3=',\n    '39+29=',\n    '6+14=',\n    '57-33=',\n    '57-29=',\n    '31-22=',\n    '44+16=',\n    '66-30=',\n    '56-55=',\n    '45+25=',\n    '1+87=',\n    '52-7=',\n    '6+28=',\n    '58-7=',\n    '65-43=',\n    '49+37=',\n    '36-2=',\n    '89-4=',\n    '66-27=',\n    '57-22=',\n    '97-3=',\n    '11+33=',\n    '47+34=',\n    '58-38=',\n    '18+71=',\n    '29+12=',\n    '11+85=',\n    '23+61=',\n    '8+24=',\n    '69-5=',\n    '67-4=',\n    '96-58=',\n    '0+81=',\n    '39-14=',\n    '78-17=',\n    '92-72=',\n    '45+43=',\n    '15-11=',\n    '17+53=',\n    '49+17=',\n    '46-8=',\n    '64-19=',\n    '69-66=',\n    '20+43=',\n    '65+24=',\n    '28+1=',\n    '69-53=',\n    '36+4=',\n    '45+10=',\n    '47-41=',\n    '72-58=',\n    '87-76=',\n    '12+8=',\n    '85+7=',\n    '47-43=',\n    '26-17=',\n    '74+8=',\n    '78-73=',\n    '22+45=',\n    '99-32=',\n    '94-17=',\n    '36+13=',\n    '53-13=',\n    '49+30=',\n    '69-20=',\n    '44+34='\n)\n\n# Word's COM Range.Text includes trailing paragraph-mark / cell-mark control\n# characters (wdParagraphEnd 0x0D, wdCellEnd 0x0D 0x07), so trim those before\n# comparing against the plain-text values pulled from the OOXML diff.\nfunction TrimMarks([string]$s) {\n    return $s.TrimEnd([char]13, [char]7)\n}\n\n# --- Update the date paragraph (first paragraph, above the table) ---\n$dateRange = $d.Paragraphs.Item(1).Range\n$dateActual = TrimMarks $dateRange.Text\nif ($dateActual -ne $dateOld) {\n    throw \"Date paragraph text mismatch: expected '$dateOld' got '$dateActual'\"\n}\n$dateRange.Text = $dateNew\n\n# --- Update each table cell, row-major (20 rows x 5 columns) ---\n$table = $d.Tables.Item(1)\n$rows = 20\n$cols = 5\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $table.Cell($r, $c)\n        $expected = $oldValues[$i]\n        $actual = TrimMarks $cell.Range.Text\n        if ($actual -ne $expected) {\n            throw \"Cell ($r,$c) text mismatch: expected '$expected' got '$actual'\"\n        }\n        $newVal = $newValues[$i]\n        if ($newVal -ne $expected) {\n            $cell.Range.Text = $newVal\n        }\n        $i++\n    }\n}\n"}
